$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data columns so numeric-looking / percentage-looking
# strings are preserved exactly as text (matching the inlineStr cells in the source).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.014.95"

# Row 3
$ws.Range("D3").Value = "2.572.96"
$ws.Range("E3").Value = "  -2.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "533.47"
$ws.Range("E5").Value = "  -0.75%  "

# Row 6
$ws.Range("D6").Value = "141.01"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  +2.76%  "

# Row 9
$ws.Range("D9").Value = "6.75"
$ws.Range("E9").Value = "  +3.17%  "

# Row 10
$ws.Range("D10").Value = "0.0989"
$ws.Range("E10").Value = "  -3.62%  "

# Row 11
$ws.Range("D11").Value = "0.138"
$ws.Range("E11").Value = "  +2.99%  "

# Row 12
$ws.Range("E12").Value = "  -1.89%  "

# Row 13
$ws.Range("D13").Value = "3.029.10"
$ws.Range("E13").Value = "  -2.50%  "

# Row 14
$ws.Range("D14").Value = "57.966.01"
$ws.Range("E14").Value = "  -2.97%  "

# Row 15
$ws.Range("D15").Value = "20.63"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16
$ws.Range("D16").Value = "2.589.78"
$ws.Range("E16").Value = "  +0.42%  "

# Row 17
$ws.Range("D17").Value = "0.0000130"
$ws.Range("E17").Value = "  -2.93%  "

# Row 18
$ws.Range("E18").Value = "  -0.72%  "

# Row 19
$ws.Range("D19").Value = "333.40"
$ws.Range("E19").Value = "  -2.64%  "

# Row 20
$ws.Range("E20").Value = "  -2.11%  "

# Row 21
$ws.Range("E21").Value = "  -4.18%  "

# Row 22
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").Value = "66.55"
$ws.Range("E23").Value = "  -1.48%  "

# Row 24
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").Value = "  +1.36%  "

# Row 25
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").Value = "0.157"
$ws.Range("E26").Value = "  -5.37%  "

# Row 27
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  -3.30%  "

# Row 28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0720"
$ws.Range("E29").Value = "  -4.17%  "

# Row 30
$ws.Range("E30").Value = "  -2.63%  "

# Row 31
$ws.Range("D31").Value = "154.76"
$ws.Range("E31").Value = "  +2.65%  "

# Row 32
$ws.Range("E32").Value = "  -0.58%  "

# Row 33
$ws.Range("E33").Value = "  -0.80%  "

# Row 34
$ws.Range("E34").Value = "  -3.99%  "

# Row 35
$ws.Range("D35").Value = "36.85"
$ws.Range("E35").Value = "  -1.59%  "

# Row 36
$ws.Range("E36").Value = "  -4.43%  "

# Row 37
$ws.Range("D37").Value = "0.825"
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("D38").Value = "0.814"
$ws.Range("E38").Value = "  -2.95%  "

# Row 39
$ws.Range("D39").Value = "1.40"
$ws.Range("E39").Value = "  -3.77%  "

# Row 40
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("D41").Value = "281.63"
$ws.Range("E41").Value = "  -2.56%  "

# Row 42
$ws.Range("E42").Value = "  +0.15%  "

# Row 43
$ws.Range("E43").Value = "  -0.80%  "

# Row 44
$ws.Range("D44").Value = "0.588"
$ws.Range("E44").Value = "  -2.58%  "

# Row 45
$ws.Range("D45").Value = "0.0944"
$ws.Range("E45").Value = "  -0.74%  "

# Row 46
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "18.23"
$ws.Range("E47").Value = "  -4.68%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0225"
$ws.Range("E48").Value = "  +0.51%  "

# Row 49
$ws.Range("D49").Value = "1.901.18"
$ws.Range("E49").Value = "  -3.30%  "

# Row 50
$ws.Range("D50").Value = "17.68"
$ws.Range("E50").Value = "  -4.35%  "

# Row 51
$ws.Range("E51").Value = "  -4.64%  "
